$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Begründung"
$ws.Range("B6").Value = "Mit Windows 10 Enterprise erhält man Einstellungen für Telemetriedaten"
$ws.Range("B6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 30

$ws.Range("B6").Select()
